$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header updates (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data updates (B2:E2) -- updated meanEMG values
$ws.Range("B2").Value = 15.225772456583631
$ws.Range("C2").Value = 10.728384784212942
$ws.Range("D2").Value = 15.426059636194196
$ws.Range("E2").Value = 8.1335694119915729

# Row 3 data updates (B3, C3 cleared, D3 new, E3 updated)
$ws.Range("B3").Value = 13.107663077411384
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 12.129492673619861
$ws.Range("E3").Value = 14.976920393192682

# Selection update to match new active range
$ws.Range("B1:E3").Select()
